$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: extend with E1:I1 = 11,12,13,14,15
$ws.Range("E1").Value = 11
$ws.Range("F1").Value = 12
$ws.Range("G1").Value = 13
$ws.Range("H1").Value = 14
$ws.Range("I1").Value = 15

# Rows 3-8: A column = 8,9,10,11,12,13
$ws.Range("A3").Value = 8
$ws.Range("A4").Value = 9
$ws.Range("A5").Value = 10
$ws.Range("A6").Value = 11
$ws.Range("A7").Value = 12
$ws.Range("A8").Value = 13
